$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6134656713070115
$ws.Cells.Item(2, 4).Value = 0.02758097233170531
$ws.Cells.Item(2, 5).Value = 0.2853314639433258
$ws.Cells.Item(2, 6).Value = 0.5018710737478926
$ws.Cells.Item(2, 7).Value = 0.3440372786528201
$ws.Cells.Item(2, 8).Value = 0.5167936572188339
$ws.Cells.Item(2, 9).Value = 0.8603903481251596
$ws.Cells.Item(2, 11).Value = 0.5121234939768442
$ws.Cells.Item(2, 12).Value = 0.1516564328683074
$ws.Cells.Item(2, 15).Value = 1.657191061463394
$ws.Cells.Item(3, 2).Value = 0.5828102936152675
$ws.Cells.Item(3, 4).Value = 0.02466433231167997
$ws.Cells.Item(3, 5).Value = 0.2885666019882382
$ws.Cells.Item(3, 6).Value = 0.5018398884598696
$ws.Cells.Item(3, 7).Value = 0.3453436757877384
$ws.Cells.Item(3, 8).Value = 0.5207963836678999
$ws.Cells.Item(3, 9).Value = 0.8723327190540857
$ws.Cells.Item(3, 11).Value = 0.4559173011872701
$ws.Cells.Item(3, 12).Value = 0.1398480638158901
$ws.Cells.Item(3, 15).Value = 1.668108969658178
$ws.Cells.Item(4, 2).Value = 0.5641571825201765
$ws.Cells.Item(4, 4).Value = 0.02286221542981792
$ws.Cells.Item(4, 5).Value = 0.2906861705798551
$ws.Cells.Item(4, 6).Value = 0.5021707890798979
$ws.Cells.Item(4, 7).Value = 0.3464275154342502
$ws.Cells.Item(4, 8).Value = 0.5234989296135382
$ws.Cells.Item(4, 9).Value = 0.8801201803366538
$ws.Cells.Item(4, 11).Value = 0.4212170410037857
$ws.Cells.Item(4, 12).Value = 0.1326422989129696
$ws.Cells.Item(4, 15).Value = 1.675914317505814
$ws.Cells.Item(5, 2).Value = 0.5565991326940889
$ws.Cells.Item(5, 4).Value = 0.02212504138914539
$ws.Cells.Item(5, 5).Value = 0.2915834357095637
$ws.Cells.Item(5, 6).Value = 0.5023936639620317
$ws.Cells.Item(5, 7).Value = 0.3469399523750738
$ws.Cells.Item(5, 8).Value = 0.5246618485142562
$ws.Cells.Item(5, 9).Value = 0.8834079986789369
$ws.Cells.Item(5, 11).Value = 0.4070297740507556
$ws.Cells.Item(5, 12).Value = 0.1297172655142305
$ws.Cells.Item(5, 15).Value = 1.679372012038428
$ws.Cells.Item(6, 2).Value = 0.5553467550995492
$ws.Cells.Item(6, 4).Value = 0.02200246651027982
$ws.Cells.Item(6, 5).Value = 0.2917344516589138
$ws.Cells.Item(6, 6).Value = 0.5024359891886263
$ws.Cells.Item(6, 7).Value = 0.347029313900336
$ws.Cells.Item(6, 8).Value = 0.5248586725051041
$ws.Cells.Item(6, 9).Value = 0.8839608458161194
$ws.Cells.Item(6, 11).Value = 0.4046712025286752
$ws.Cells.Item(6, 12).Value = 0.1292322574589377
$ws.Cells.Item(6, 15).Value = 1.679962884385347
$ws.Cells.Item(7, 2).Value = 0.5640550760678593
$ws.Cells.Item(7, 4).Value = 0.02285228491260938
$ws.Cells.Item(7, 5).Value = 0.2906981356287517
$ws.Cells.Item(7, 6).Value = 0.5021734384078016
$ws.Cells.Item(7, 7).Value = 0.3464341399183368
$ws.Cells.Item(7, 8).Value = 0.5235143636389807
$ws.Cells.Item(7, 9).Value = 0.8801640579430146
$ws.Cells.Item(7, 11).Value = 0.4210258940260871
$ws.Cells.Item(7, 12).Value = 0.1326028046312331
$ws.Cells.Item(7, 15).Value = 1.675959827925112
$ws.Cells.Item(8, 2).Value = 0.6028609104277791
$ws.Cells.Item(8, 4).Value = 0.02657768190180576
$ws.Cells.Item(8, 5).Value = 0.2864193187389983
$ws.Cells.Item(8, 6).Value = 0.501787671628918
$ws.Cells.Item(8, 7).Value = 0.3444292111256289
$ws.Cells.Item(8, 8).Value = 0.5181230005377344
$ws.Cells.Item(8, 9).Value = 0.8644136869760182
$ws.Cells.Item(8, 11).Value = 0.4927835001551841
$ws.Cells.Item(8, 12).Value = 0.1475757345735929
$ws.Cells.Item(8, 15).Value = 1.660726846081047
$ws.Cells.Item(9, 2).Value = 0.6802791148505776
$ws.Cells.Item(9, 4).Value = 0.03379209729747146
$ws.Cells.Item(9, 5).Value = 0.2790837084437374
$ws.Cells.Item(9, 6).Value = 0.5038091260055211
$ws.Cells.Item(9, 7).Value = 0.3427365136864395
$ws.Cells.Item(9, 8).Value = 0.5094917723296035
$ws.Cells.Item(9, 9).Value = 0.8371354069397796
$ws.Cells.Item(9, 11).Value = 0.6319592374143781
$ws.Cells.Item(9, 12).Value = 0.1772864239484306
$ws.Cells.Item(9, 15).Value = 1.63960141460214
$ws.Cells.Item(10, 2).Value = 0.7379366500385061
$ws.Cells.Item(10, 4).Value = 0.03903544905944045
$ws.Cells.Item(10, 5).Value = 0.2743352091649953
$ws.Cells.Item(10, 6).Value = 0.5069896059651171
$ws.Cells.Item(10, 7).Value = 0.3428636008088191
$ws.Cells.Item(10, 8).Value = 0.5043318092783906
$ws.Cells.Item(10, 9).Value = 0.8192925179405233
$ws.Cells.Item(10, 11).Value = 0.7332312997810106
$ws.Cells.Item(10, 12).Value = 0.1993227400539865
$ws.Cells.Item(10, 15).Value = 1.629420978928877
$ws.Cells.Item(11, 2).Value = 0.7643303633581127
$ws.Cells.Item(11, 4).Value = 0.04140808790062067
$ws.Cells.Item(11, 5).Value = 0.272313691262541
$ws.Cells.Item(11, 6).Value = 0.5088050568488214
$ws.Cells.Item(11, 7).Value = 0.3432203178360993
$ws.Cells.Item(11, 8).Value = 0.5022405251923772
$ws.Cells.Item(11, 9).Value = 0.8116524382664512
$ws.Cells.Item(11, 11).Value = 0.7790816136239584
$ws.Cells.Item(11, 12).Value = 0.2093918867037985
$ws.Cells.Item(11, 15).Value = 1.625951146910808
$ws.Cells.Item(12, 2).Value = 0.774348134002139
$ws.Cells.Item(12, 4).Value = 0.04230469750500276
$ws.Cells.Item(12, 5).Value = 0.2715680904128863
$ws.Cells.Item(12, 6).Value = 0.5095455411002945
$ws.Cells.Item(12, 7).Value = 0.3433984684639881
$ws.Cells.Item(12, 8).Value = 0.5014853906840955
$ws.Cells.Item(12, 9).Value = 0.8088278925176997
$ws.Cells.Item(12, 11).Value = 0.7964115738862461
$ws.Cells.Item(12, 12).Value = 0.2132111139191295
$ws.Cells.Item(12, 15).Value = 1.624804329712219
$ws.Cells.Item(13, 2).Value = 0.7721896138346551
$ws.Cells.Item(13, 4).Value = 0.04211167993479137
$ws.Cells.Item(13, 5).Value = 0.2717277840719436
$ws.Cells.Item(13, 6).Value = 0.5093837067936562
$ws.Cells.Item(13, 7).Value = 0.3433581836271031
$ws.Cells.Item(13, 8).Value = 0.5016463865836869
$ws.Cells.Item(13, 9).Value = 0.8094331580664313
$ws.Cells.Item(13, 11).Value = 0.7926807221865317
$ws.Cells.Item(13, 12).Value = 0.2123882993061414
$ws.Cells.Item(13, 15).Value = 1.625043881850502
$ws.Cells.Item(14, 2).Value = 0.7651540734523508
$ws.Cells.Item(14, 4).Value = 0.0414818900035101
$ws.Cells.Item(14, 5).Value = 0.2722519514734945
$ws.Cells.Item(14, 6).Value = 0.5088649144738682
$ws.Cells.Item(14, 7).Value = 0.343234110719699
$ws.Cells.Item(14, 8).Value = 0.5021776626479664
$ws.Cells.Item(14, 9).Value = 0.8114186864269683
$ws.Cells.Item(14, 11).Value = 0.7805080186183204
$ws.Cells.Item(14, 12).Value = 0.2097059723433432
$ws.Cells.Item(14, 15).Value = 1.625853447190678
$ws.Cells.Item(15, 2).Value = 0.7608475825758774
$ws.Cells.Item(15, 4).Value = 0.04109588222212324
$ws.Cells.Item(15, 5).Value = 0.2725756104691364
$ws.Cells.Item(15, 6).Value = 0.5085540430447537
$ws.Cells.Item(15, 7).Value = 0.3431637239920846
$ws.Cells.Item(15, 8).Value = 0.5025078749166241
$ws.Cells.Item(15, 9).Value = 0.8126438128280906
$ws.Cells.Item(15, 11).Value = 0.7730476139269058
$ws.Cells.Item(15, 12).Value = 0.2080637797961629
$ws.Cells.Item(15, 15).Value = 1.626371099054168
$ws.Cells.Item(16, 2).Value = 0.7362150207314926
$ws.Cells.Item(16, 4).Value = 0.03888013420322523
$ws.Cells.Item(16, 5).Value = 0.2744701068284776
$ws.Cells.Item(16, 6).Value = 0.5068783801612753
$ws.Cells.Item(16, 7).Value = 0.3428463107497421
$ws.Cells.Item(16, 8).Value = 0.5044736290847496
$ws.Cells.Item(16, 9).Value = 0.8198014110785135
$ws.Cells.Item(16, 11).Value = 0.7302303769605771
$ws.Cells.Item(16, 12).Value = 0.198665582898272
$ws.Cells.Item(16, 15).Value = 1.629671120919312
$ws.Cells.Item(17, 2).Value = 0.7211455179944153
$ws.Cells.Item(17, 4).Value = 0.03751758695523222
$ws.Cells.Item(17, 5).Value = 0.2756677986924885
$ws.Cells.Item(17, 6).Value = 0.5059448342115687
$ws.Cells.Item(17, 7).Value = 0.342728207349758
$ws.Cells.Item(17, 8).Value = 0.5057451083164963
$ws.Cells.Item(17, 9).Value = 0.8243144931037296
$ws.Cells.Item(17, 11).Value = 0.7039065597549268
$ws.Cells.Item(17, 12).Value = 0.1929114253477877
$ws.Cells.Item(17, 15).Value = 1.631993114325994
$ws.Cells.Item(18, 2).Value = 0.7124935110565502
$ws.Cells.Item(18, 4).Value = 0.03673270377359472
$ws.Cells.Item(18, 5).Value = 0.2763697264587694
$ws.Cells.Item(18, 6).Value = 0.5054425808178564
$ws.Cells.Item(18, 7).Value = 0.3426884065583451
$ws.Cells.Item(18, 8).Value = 0.5065005277831176
$ws.Cells.Item(18, 9).Value = 0.8269551696396888
$ws.Cells.Item(18, 11).Value = 0.6887452613599407
$ws.Cells.Item(18, 12).Value = 0.1896060051869881
$ws.Cells.Item(18, 15).Value = 1.633437962370039
$ws.Cells.Item(19, 2).Value = 0.709566785625185
$ws.Cells.Item(19, 4).Value = 0.03646675423100021
$ws.Cells.Item(19, 5).Value = 0.2766096287861721
$ws.Cells.Item(19, 6).Value = 0.5052784856617691
$ws.Cells.Item(19, 7).Value = 0.3426797592141355
$ws.Cells.Item(19, 8).Value = 0.5067604394727212
$ws.Cells.Item(19, 9).Value = 0.8278569630248906
$ws.Cells.Item(19, 11).Value = 0.6836084125338289
$ws.Cells.Item(19, 12).Value = 0.1884875765248637
$ws.Cells.Item(19, 15).Value = 1.633945930973155
$ws.Cells.Item(20, 2).Value = 0.722748085556276
$ws.Cells.Item(20, 4).Value = 0.0376627551128621
$ws.Cells.Item(20, 5).Value = 0.2755389522146832
$ws.Cells.Item(20, 6).Value = 0.5060406206777017
$ws.Cells.Item(20, 7).Value = 0.3427378677569379
$ws.Cells.Item(20, 8).Value = 0.50560726328046
$ws.Cells.Item(20, 9).Value = 0.8238294235932457
$ws.Cells.Item(20, 11).Value = 0.7067109082182128
$ws.Cells.Item(20, 12).Value = 0.1935235298161047
$ws.Cells.Item(20, 15).Value = 1.631734621121808
$ws.Cells.Item(21, 2).Value = 0.7672199622388121
$ws.Cells.Item(21, 4).Value = 0.04166692535292782
$ws.Cells.Item(21, 5).Value = 0.2720974508181904
$ws.Cells.Item(21, 6).Value = 0.5090158575951946
$ws.Cells.Item(21, 7).Value = 0.3432693843524319
$ws.Cells.Item(21, 8).Value = 0.5020206157311904
$ws.Cells.Item(21, 9).Value = 0.8108336269422303
$ws.Cells.Item(21, 11).Value = 0.784084328999711
$ws.Cells.Item(21, 12).Value = 0.2104936682590903
$ws.Cells.Item(21, 15).Value = 1.625611121192833
$ws.Cells.Item(22, 2).Value = 0.7964189073159957
$ws.Cells.Item(22, 4).Value = 0.04427303482749778
$ws.Cells.Item(22, 5).Value = 0.2699642349779516
$ws.Cells.Item(22, 6).Value = 0.511269344085882
$ws.Cells.Item(22, 7).Value = 0.3438678435366285
$ws.Cells.Item(22, 8).Value = 0.4998909595336016
$ws.Cells.Item(22, 9).Value = 0.8027399126161416
$ws.Cells.Item(22, 11).Value = 0.8344620548372745
$ws.Cells.Item(22, 12).Value = 0.2216210430476195
$ws.Cells.Item(22, 15).Value = 1.622583300688916
$ws.Cells.Item(23, 2).Value = 0.780822845254022
$ws.Cells.Item(23, 4).Value = 0.04288311324378924
$ws.Cells.Item(23, 5).Value = 0.2710921666474935
$ws.Cells.Item(23, 6).Value = 0.5100383420408789
$ws.Cells.Item(23, 7).Value = 0.343525432007894
$ws.Cells.Item(23, 8).Value = 0.5010079852833798
$ws.Cells.Item(23, 9).Value = 0.8070230889145567
$ws.Cells.Item(23, 11).Value = 0.8075922830063007
$ws.Cells.Item(23, 12).Value = 0.2156788786654289
$ws.Cells.Item(23, 15).Value = 1.624110117173075
$ws.Cells.Item(24, 2).Value = 0.7220235284191006
$ws.Cells.Item(24, 4).Value = 0.03759712937147697
$ws.Cells.Item(24, 5).Value = 0.2755971620946998
$ws.Cells.Item(24, 6).Value = 0.5059972082939552
$ws.Cells.Item(24, 7).Value = 0.3427334127693413
$ws.Cells.Item(24, 8).Value = 0.5056695069247752
$ws.Cells.Item(24, 9).Value = 0.8240485800936614
$ws.Cells.Item(24, 11).Value = 0.7054431474976752
$ws.Cells.Item(24, 12).Value = 0.1932467888964027
$ws.Cells.Item(24, 15).Value = 1.631851143570799
$ws.Cells.Item(25, 2).Value = 0.6591968482607342
$ws.Cells.Item(25, 4).Value = 0.03185031791193182
$ws.Cells.Item(25, 5).Value = 0.2809554602489488
$ws.Cells.Item(25, 6).Value = 0.5029646787741555
$ws.Cells.Item(25, 7).Value = 0.3429541345888509
$ws.Cells.Item(25, 8).Value = 0.5116191359966535
$ws.Cells.Item(25, 9).Value = 0.8441287474948567
$ws.Cells.Item(25, 11).Value = 0.5944777219969524
$ws.Cells.Item(25, 12).Value = 0.1692120446508909
$ws.Cells.Item(25, 15).Value = 1.644379164620034
